$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Update "Sr. No." column (A2:A5) from dates/odd numbers to sequential 1..4
# A2 previously carried a date-number style; clear it by copying the
# (unstyled) format from A3 before writing the new plain number.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# Update Test Case Title column (F2:F5)
$ws.Range("F2").Value = "Title 1"
$ws.Range("F3").Value = "Title 2"
$ws.Range("F5").Value = "Title 4"
$ws.Range("F4").Value = "Smoke"

# Update the active selection to F4
$ws.Range("F4").Select()
